$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H28").Value = 1809.25
$ws.Range("I28").Value = 1924.2858
$ws.Range("J28").Value = 1004
$ws.Range("K28").Value = 1924.2858
$ws.Range("L28").Value = 1004
$ws.Range("M28").Value = -1439.2858
$ws.Range("N28").Value = -1974
$ws.Range("H98").Value = 1879.5217
$ws.Range("I98").Value = 1437.8422
$ws.Range("K98").Value = 1437.8422
$ws.Range("M98").Value = 60.15779999999995
$ws.Range("H122").Value = 1879.5217
$ws.Range("I122").Value = 1437.8422
$ws.Range("K122").Value = 4313.5266
$ws.Range("M122").Value = -1863.5266
$ws.Range("H137").Value = 1793973.4
$ws.Range("I137").Value = 1673.2727
$ws.Range("J137").Value = 6175151
$ws.Range("K137").Value = 5019.8181
$ws.Range("L137").Value = 18525453
$ws.Range("M137").Value = -2469.8181
$ws.Range("N137").Value = -18530553

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H2").Value = 3430.1428
$ws.Range("I2").Value = 3622.2
$ws.Range("J2").Value = 2950
$ws.Range("K2").Value = 3622.2
$ws.Range("L2").Value = 2950
$ws.Range("M2").Value = -3509.2
$ws.Range("N2").Value = -3176
$ws.Range("H32").Value = 3052.5254
$ws.Range("I32").Value = 3132.027
$ws.Range("J32").Value = 2918.818
$ws.Range("K32").Value = 3132.027
$ws.Range("L32").Value = 2918.818
$ws.Range("M32").Value = -2845.027
$ws.Range("N32").Value = -3492.818
$ws.Range("H116").Value = 3430.1428
$ws.Range("I116").Value = 3622.2
$ws.Range("J116").Value = 2950
$ws.Range("K116").Value = 3622.2
$ws.Range("L116").Value = 2950
$ws.Range("M116").Value = -1328.2
$ws.Range("N116").Value = -7538

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H3").Value = 3430.1428
$ws.Range("I3").Value = 3622.2
$ws.Range("J3").Value = 2950
$ws.Range("K3").Value = 3622.2
$ws.Range("L3").Value = 2950
$ws.Range("M3").Value = -3508.2
$ws.Range("N3").Value = -3178
$ws.Range("H20").Value = 1439.25
$ws.Range("I20").Value = 1191.5
$ws.Range("J20").Value = 1604.4166
$ws.Range("K20").Value = 1191.5
$ws.Range("L20").Value = 1604.4166
$ws.Range("M20").Value = -944.5
$ws.Range("N20").Value = -2098.4166
$ws.Range("H21").Value = 21500
$ws.Range("J21").Value = 21500
$ws.Range("L21").Value = 21500
$ws.Range("N21").Value = -21972
$ws.Range("H40").Value = 27499.75
$ws.Range("J40").Value = 27499.75
$ws.Range("L40").Value = 27499.75
$ws.Range("N40").Value = -28029.75
$ws.Range("H97").Value = 10405.4
$ws.Range("I97").Value = 2864.8572
$ws.Range("K97").Value = 2864.8572
$ws.Range("M97").Value = -1873.8572
$ws.Range("H137").Value = 35940.2
$ws.Range("J137").Value = 35940.2
$ws.Range("L137").Value = 35940.2
$ws.Range("N137").Value = -46140.2

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H16").Value = 1048.7778
$ws.Range("I16").Value = 1061.5
$ws.Range("J16").Value = 1038.6
$ws.Range("K16").Value = 1061.5
$ws.Range("L16").Value = 1038.6
$ws.Range("M16").Value = -774.5
$ws.Range("N16").Value = -1612.6
$ws.Range("H31").Value = 1584.5161
$ws.Range("I31").Value = 1293.5
$ws.Range("J31").Value = 2295.889
$ws.Range("K31").Value = 1293.5
$ws.Range("L31").Value = 2295.889
$ws.Range("M31").Value = -998.5
$ws.Range("N31").Value = -2885.889
$ws.Range("H34").Value = 1584.5161
$ws.Range("I34").Value = 1293.5
$ws.Range("J34").Value = 2295.889
$ws.Range("K34").Value = 1293.5
$ws.Range("L34").Value = 2295.889
$ws.Range("M34").Value = -1091.5
$ws.Range("N34").Value = -2699.889
$ws.Range("H58").Value = 1668.55
$ws.Range("I58").Value = 1679.625
$ws.Range("J58").Value = 1624.25
$ws.Range("K58").Value = 1679.625
$ws.Range("L58").Value = 1624.25
$ws.Range("M58").Value = -1476.625
$ws.Range("N58").Value = -2030.25
$ws.Range("H113").Value = 1048.7778
$ws.Range("I113").Value = 1061.5
$ws.Range("J113").Value = 1038.6
$ws.Range("K113").Value = 1061.5
$ws.Range("L113").Value = 1038.6
$ws.Range("M113").Value = 1108.5
$ws.Range("N113").Value = -5378.6
$ws.Range("H122").Value = 994.15
$ws.Range("I122").Value = 945.5333000000001
$ws.Range("J122").Value = 1140
$ws.Range("K122").Value = 2836.5999
$ws.Range("L122").Value = 3420
$ws.Range("M122").Value = -386.5999000000002
$ws.Range("N122").Value = -8320
$ws.Range("H132").Value = 3375.375
$ws.Range("I132").Value = 3069.2307
$ws.Range("J132").Value = 4702
$ws.Range("K132").Value = 9207.6921
$ws.Range("L132").Value = 14106
$ws.Range("M132").Value = -6677.6921
$ws.Range("N132").Value = -19166
$ws.Range("H134").Value = 1741.1714
$ws.Range("I134").Value = 1671.5667
$ws.Range("J134").Value = 2158.8
$ws.Range("K134").Value = 5014.7001
$ws.Range("L134").Value = 6476.400000000001
$ws.Range("M134").Value = -2479.7001
$ws.Range("N134").Value = -11546.4
$ws.Range("H136").Value = 1668.55
$ws.Range("I136").Value = 1679.625
$ws.Range("J136").Value = 1624.25
$ws.Range("K136").Value = 5038.875
$ws.Range("L136").Value = 4872.75
$ws.Range("M136").Value = -2488.875
$ws.Range("N136").Value = -9972.75

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H131").Value = 881.4299999999999
$ws.Range("J131").Value = 897.49475
$ws.Range("L131").Value = 2692.48425
$ws.Range("N131").Value = -12772.48425

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H80").Value = 2810.2334
$ws.Range("I80").Value = 2699.7827
$ws.Range("J80").Value = 3173.1428
$ws.Range("K80").Value = 2699.7827
$ws.Range("L80").Value = 3173.1428
$ws.Range("M80").Value = -1701.7827
$ws.Range("N80").Value = -5169.1428
$ws.Range("H83").Value = 2810.2334
$ws.Range("I83").Value = 2699.7827
$ws.Range("J83").Value = 3173.1428
$ws.Range("K83").Value = 13498.9135
$ws.Range("L83").Value = 15865.714
$ws.Range("M83").Value = -8506.913500000001
$ws.Range("N83").Value = -25849.714

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H40").Value = 5911.5557
$ws.Range("I40").Value = 5911.5557
$ws.Range("K40").Value = 5911.5557
$ws.Range("M40").Value = -5775.5557
$ws.Range("H106").Value = 29500
$ws.Range("J106").Value = 29500
$ws.Range("L106").Value = 29500
$ws.Range("N106").Value = -32024
$ws.Range("H122").Value = 3538.5
$ws.Range("I122").Value = 2742.8572
$ws.Range("J122").Value = 4044.818
$ws.Range("K122").Value = 8228.571599999999
$ws.Range("L122").Value = 12134.454
$ws.Range("M122").Value = -5778.571599999999
$ws.Range("N122").Value = -17034.454

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H101").Value = 29999.666
$ws.Range("J101").Value = 29999.666
$ws.Range("L101").Value = 29999.666
$ws.Range("N101").Value = -36489.666
$ws.Range("H122").Value = 172416.58
$ws.Range("I122").Value = 1382.2
$ws.Range("J122").Value = 600002.5
$ws.Range("K122").Value = 4146.6
$ws.Range("L122").Value = 1800007.5
$ws.Range("M122").Value = -1696.6
$ws.Range("N122").Value = -1804907.5
